$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controllers")

# Update the PVSystem / device identifiers and the sample selector after
# finishing Volt-Var debugging on a different device/index.
$ws.Range("B3").Value = "PVSystem.oh_261584_2_4"
$ws.Range("B4").Value = "dev_261585_2_4"
$ws.Range("F4").Value = "['Even']"

# Move the active selection to F4 on the Controllers sheet.
$ws.Activate()
$ws.Range("F4").Select()
